# Auto-generated: update leve profit calc columns (H-N) per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (item id 5487)
$ws.Range("H9").Value = 58.5
$ws.Range("I9").Value = 58.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 58.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 110.5
$ws.Range("N9").Value = $null
# Row 33 (item id 5512)
$ws.Range("H33").Value = 117.375
$ws.Range("J33").Value = 200
$ws.Range("L33").Value = 200
$ws.Range("N33").Value = -658
# Row 105 (item id 18668)
$ws.Range("H105").Value = 25509.8
$ws.Range("J105").Value = 25509.8
$ws.Range("L105").Value = 25509.8
$ws.Range("N105").Value = -32497.8
# Row 107 (item id 1454.25)
$ws.Range("I107").Value = 1722.7
$ws.Range("J107").Value = 112
$ws.Range("K107").Value = 1722.7
$ws.Range("L107").Value = 112
$ws.Range("M107").Value = 197.3
$ws.Range("N107").Value = -3952
# Row 132 (item id 44049)
$ws.Range("H132").Value = 3415.9524
$ws.Range("I132").Value = 1364.2307
$ws.Range("K132").Value = 4092.6921
$ws.Range("M132").Value = -1562.6921
# Row 138 (item id 44169)
$ws.Range("H138").Value = 3173
$ws.Range("J138").Value = 5043.4
$ws.Range("L138").Value = 15130.2
$ws.Range("N138").Value = -25410.2

$ws = $wb.Worksheets.Item("ARM")
# Row 35 (item id 2473)
$ws.Range("H35").Value = 12884.25
$ws.Range("I35").Value = 12884.25
$ws.Range("K35").Value = 12884.25
$ws.Range("M35").Value = -12478.25
# Row 74 (item id 44000)
$ws.Range("H74").Value = 8013.6665
$ws.Range("I74").Value = 8013.6665
$ws.Range("K74").Value = 8013.6665
$ws.Range("M74").Value = -7139.6665
# Row 77 (item id 44000)
$ws.Range("H77").Value = 8013.6665
$ws.Range("I77").Value = 8013.6665
$ws.Range("K77").Value = 40068.3325
$ws.Range("M77").Value = -35700.3325
# Row 88 (item id 12530)
$ws.Range("H88").Value = 1062.6666
$ws.Range("I88").Value = 1046.2858
$ws.Range("K88").Value = 1046.2858
$ws.Range("M88").Value = -640.2858000000001
# Row 91 (item id 12530)
$ws.Range("H91").Value = 1062.6666
$ws.Range("I91").Value = 1046.2858
$ws.Range("K91").Value = 1046.2858
$ws.Range("M91").Value = 357.7141999999999
# Row 97 (item id 19941)
$ws.Range("H97").Value = 1111
$ws.Range("I97").Value = 1111
$ws.Range("K97").Value = 1111
$ws.Range("M97").Value = -615
# Row 132 (item id 43997)
$ws.Range("H132").Value = 1685.9333
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
# Row 41 (item id 22899)
$ws.Range("H41").Value = 240000
$ws.Range("J41").Value = 240000
$ws.Range("L41").Value = 240000
$ws.Range("N41").Value = -240776
# Row 75 (item id 11872)
$ws.Range("H75").Value = 3266.3333
$ws.Range("I75").Value = 3266.3333
$ws.Range("K75").Value = 3266.3333
$ws.Range("M75").Value = -2330.3333
# Row 78 (item id 11872)
$ws.Range("H78").Value = 3266.3333
$ws.Range("I78").Value = 3266.3333
$ws.Range("K78").Value = 9798.999899999999
$ws.Range("M78").Value = -5118.999899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (item id 27691)
$ws.Range("H16").Value = 1340.6428
$ws.Range("I16").Value = 1367.6364
$ws.Range("J16").Value = 1241.6666
$ws.Range("K16").Value = 1367.6364
$ws.Range("L16").Value = 1241.6666
$ws.Range("M16").Value = -1080.6364
$ws.Range("N16").Value = -1815.6666
# Row 22 (item id 5367)
$ws.Range("H22").Value = 2554.6667
$ws.Range("I22").Value = 2499.6
$ws.Range("J22").Value = 2623.5
$ws.Range("K22").Value = 2499.6
$ws.Range("L22").Value = 2623.5
$ws.Range("M22").Value = -2149.6
$ws.Range("N22").Value = -3323.5
# Row 39 (item id 1915)
$ws.Range("H39").Value = 6037
$ws.Range("I39").Value = 425.5
$ws.Range("K39").Value = 425.5
$ws.Range("M39").Value = -34.5
# Row 49 (item id 1915)
$ws.Range("H49").Value = 6037
$ws.Range("I49").Value = 425.5
$ws.Range("K49").Value = 425.5
$ws.Range("M49").Value = -243.5
# Row 58 (item id 44021)
$ws.Range("H58").Value = 8995
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null
# Row 62 (item id 12580)
$ws.Range("H62").Value = 9083.833000000001
$ws.Range("I62").Value = 8999.5
$ws.Range("K62").Value = 8999.5
$ws.Range("M62").Value = -8375.5
# Row 65 (item id 12580)
$ws.Range("H65").Value = 9083.833000000001
$ws.Range("I65").Value = 8999.5
$ws.Range("K65").Value = 44997.5
$ws.Range("M65").Value = -41877.5
# Row 105 (item id 19928)
$ws.Range("H105").Value = 1111
$ws.Range("I105").Value = 1111
$ws.Range("K105").Value = 1111
$ws.Range("M105").Value = 636
# Row 113 (item id 27691)
$ws.Range("H113").Value = 1340.6428
$ws.Range("I113").Value = 1367.6364
$ws.Range("J113").Value = 1241.6666
$ws.Range("K113").Value = 1367.6364
$ws.Range("L113").Value = 1241.6666
$ws.Range("M113").Value = 802.3635999999999
$ws.Range("N113").Value = -5581.6666
# Row 136 (item id 44021)
$ws.Range("H136").Value = 8995
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 35 (item id 4718)
$ws.Range("H35").Value = 1000000
$ws.Range("J35").Value = 1000000
$ws.Range("L35").Value = 3000000
$ws.Range("N35").Value = -3000576
# Row 116 (item id 27866)
$ws.Range("H116").Value = 2666.3333
$ws.Range("I116").Value = 2666.3333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 7998.999899999999
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -4556.999899999999
$ws.Range("N116").Value = $null
# Row 141 (item id 44076)
$ws.Range("H141").Value = 8250
$ws.Range("I141").Value = 1500
$ws.Range("K141").Value = 4500
$ws.Range("M141").Value = 680

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (item id 19940)
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("N97").Value = $null
# Row 101 (item id 18513)
$ws.Range("H101").Value = 18995
$ws.Range("J101").Value = 18995
$ws.Range("L101").Value = 18995
$ws.Range("N101").Value = -25485
# Row 113 (item id 27710)
$ws.Range("H113").Value = 5244.125
$ws.Range("I113").Value = 4446.7
$ws.Range("K113").Value = 4446.7
$ws.Range("M113").Value = -2276.7

$ws = $wb.Worksheets.Item("LTW")
# Row 45 (item id 3851)
$ws.Range("H45").Value = 4041
$ws.Range("I45").Value = 4041
$ws.Range("K45").Value = 4041
$ws.Range("M45").Value = -3634
# Row 46 (item id 5282)
$ws.Range("H46").Value = 999.5
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811
# Row 93 (item id 19993)
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").Value = $null
# Row 106 (item id 18713)
$ws.Range("H106").Value = 9667.25
$ws.Range("J106").Value = 9667.25
$ws.Range("L106").Value = 9667.25
$ws.Range("N106").Value = -12191.25
# Row 122 (item id 36247)
$ws.Range("H122").Value = 2998.5
$ws.Range("I122").Value = 2998.5
$ws.Range("K122").Value = 8995.5
$ws.Range("M122").Value = -6545.5

$ws = $wb.Worksheets.Item("WVR")
# Row 32 (item id 3066)
$ws.Range("H32").Value = 9000
$ws.Range("I32").Value = 9000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -8683
$ws.Range("N32").Value = $null
# Row 96 (item id 19977)
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 4000
$ws.Range("K96").Value = 4000
$ws.Range("M96").Value = -2627
# Row 104 (item id 18691)
$ws.Range("H104").Value = 31584.834
$ws.Range("J104").Value = 31584.834
$ws.Range("L104").Value = 31584.834
$ws.Range("N104").Value = -38572.834
# Row 105 (item id 18710)
$ws.Range("H105").Value = 12666.667
$ws.Range("J105").Value = 12666.667
$ws.Range("L105").Value = 12666.667
$ws.Range("N105").Value = -19654.667
# Row 122 (item id 36208)
$ws.Range("H122").Value = 4668
$ws.Range("I122").Value = 4668
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14004
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11554
$ws.Range("N122").Value = $null
# Row 126 (item id 36210)
$ws.Range("H126").Value = 5356.3887
$ws.Range("I126").Value = 3601.7778
$ws.Range("J126").Value = 7111
$ws.Range("K126").Value = 10805.3334
$ws.Range("L126").Value = 21333
$ws.Range("M126").Value = -8335.3334
$ws.Range("N126").Value = -26273
# Row 132 (item id 44029)
$ws.Range("H132").Value = 3002.5
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
